$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Fitness") values for rows 2-57 (Generation 0-55) change to 7295
for ($r = 2; $r -le 57; $r++) {
    $ws.Cells.Item($r, 3).Value = 7295
}

# Column C ("Fitness") values for rows 58-252 (Generation 56-250) change to 7293
for ($r = 58; $r -le 252; $r++) {
    $ws.Cells.Item($r, 3).Value = 7293
}
